$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Add a new LeetCode entry (row 20 on Sheet1) ---
$ws1.Range("B20").Value = "Binary Tree Inorder Traversal"
$ws1.Range("C20").Value = 1
$ws1.Range("D20").Value = 1
$ws1.Range("E20").Value = 37
$ws1.Range("F20").Value = 0.33
$ws1.Range("G20").Value = 16.1
$ws1.Range("H20").Value = 0.0125
$ws1.Range("I20").Value = "https://leetcode.com/problems/binary-tree-inorder-traversal/submissions/1060854991/"
$ws1.Range("J20").Value = "stack, recursion, depth-first search"

# --- Extend the summary averages on Sheet2 to cover the new row ---
$ws2.Range("B3").Formula = "=SUM(Sheet1!C3:C500)/COUNT(Sheet1!C3:C500)"
$ws2.Range("C3").Formula = "=SUM(Sheet1!D3:D500)/COUNT(Sheet1!D3:D500)"
$ws2.Range("D3").Formula = "=SUM(Sheet1!E3:E500)/COUNT(Sheet1!E3:E500)"
$ws2.Range("E3").Formula = "=SUM(Sheet1!F3:F500)/COUNT(Sheet1!F3:F500)"
$ws2.Range("F3").Formula = "=SUM(Sheet1!G3:G500)/COUNT(Sheet1!G3:G500)"
$ws2.Range("G3").Formula = "=SUM(Sheet1!H3:H500)/COUNT(Sheet1!H3:H500)"

# --- Restore the last-saved cursor position on Sheet1, then switch the
#     active tab/selection to Sheet2 (matches the saved workbook view state) ---
$ws1.Activate()
$ws1.Range("D24").Select()

$ws2.Activate()
$ws2.Range("G6").Select()
